{"js": "// Replace the date line and each two-digit multiplication problem in the\n// table with the new values from the commit. Every search string below is\n// unique within the document, so a plain exact-match search/replace is\n// unambiguous and keeps the original run formatting (font/size) intact.\nconst replacements = [\n  [\"2024-01-19 Friday\", \"2024-01-20 Saturday\"],\n  [\"90\u00d738=\", \"71\u00d724=\"],\n  [\"93\u00d721=\", \"60\u00d763=\"],\n  [\"83\u00d724=\", \"82\u00d754=\"],\n  [\"63\u00d781=\", \"19\u00d722=\"],\n  [\"13\u00d772=\", \"71\u00d762=\"],\n  [\"64\u00d721=\", \"26\u00d765=\"],\n  [\"89\u00d742=\", \"86\u00d780=\"],\n  [\"83\u00d743=\", \"87\u00d788=\"],\n  [\"70\u00d797=\", \"44\u00d713=\"],\n  [\"65\u00d744=\", \"77\u00d760=\"],\n  [\"74\u00d724=\", \"63\u00d715=\"],\n  [\"28\u00d775=\", \"64\u00d757=\"],\n  [\"32\u00d784=\", \"68\u00d711=\"],\n  [\"15\u00d798=\", \"24\u00d791=\"],\n  [\"31\u00d723=\", \"22\u00d781=\"],\n  [\"81\u00d768=\", \"42\u00d727=\"],\n  [\"49\u00d767=\", \"36\u00d786=\"],\n  [\"24\u00d787=\", \"45\u00d726=\"],\n  [\"39\u00d781=\", \"64\u00d733=\"],\n  [\"89\u00d755=\", \"21\u00d754=\"],\n  [\"89\u00d720=\", \"17\u00d742=\"],\n  [\"32\u00d760=\", \"62\u00d768=\"],\n  [\"84\u00d730=\", \"24\u00d745=\"],\n  [\"91\u00d721=\", \"86\u00d745=\"],\n  [\"66\u00d798=\", \"92\u00d713=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit multiplication problem in the\n# table to the new values from the commit. Each old value is unique in the\n# document, so Find/Replace (ReplaceAll) is unambiguous and preserves the\n# existing run formatting (font/size) of the matched text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-01-19 Friday\"; New = \"2024-01-20 Saturday\" },\n    @{ Old = \"90\u00d738=\"; New = \"71\u00d724=\" },\n    @{ Old = \"93\u00d721=\"; New = \"60\u00d763=\" },\n    @{ Old = \"83\u00d724=\"; New = \"82\u00d754=\" },\n    @{ Old = \"63\u00d781=\"; New = \"19\u00d722=\" },\n    @{ Old = \"13\u00d772=\"; New = \"71\u00d762=\" },\n    @{ Old = \"64\u00d721=\"; New = \"26\u00d765=\" },\n    @{ Old = \"89\u00d742=\"; New = \"86\u00d780=\" },\n    @{ Old = \"83\u00d743=\"; New = \"87\u00d788=\" },\n    @{ Old = \"70\u00d797=\"; New = \"44\u00d713=\" },\n    @{ Old = \"65\u00d744=\"; New = \"77\u00d760=\" },\n    @{ Old = \"74\u00d724=\"; New = \"63\u00d715=\" },\n    @{ Old = \"28\u00d775=\"; New = \"64\u00d757=\" },\n    @{ Old = \"32\u00d784=\"; New = \"68\u00d711=\" },\n    @{ Old = \"15\u00d798=\"; New = \"24\u00d791=\" },\n    @{ Old = \"31\u00d723=\"; New = \"22\u00d781=\" },\n    @{ Old = \"81\u00d768=\"; New = \"42\u00d727=\" },\n    @{ Old = \"49\u00d767=\"; New = \"36\u00d786=\" },\n    @{ Old = \"24\u00d787=\"; New = \"45\u00d726=\" },\n    @{ Old = \"39\u00d781=\"; New = \"64\u00d733=\" },\n    @{ Old = \"89\u00d755=\"; New = \"21\u00d754=\" },\n    @{ Old = \"89\u00d720=\"; New = \"17\u00d742=\" },\n    @{ Old = \"32\u00d760=\"; New = \"62\u00d768=\" },\n    @{ Old = \"84\u00d730=\"; New = \"24\u00d745=\" },\n    @{ Old = \"91\u00d721=\"; New = \"86\u00d745=\" },\n    @{ Old = \"66\u00d798=\"; New = \"92\u00d713=\" }\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$wdFindContinue, [ref]$false, [ref]$find.Replacement.Text, [ref]$wdReplaceAll) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
